$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 16:52"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 889391
$ws.Range("C4").Value = 2949
$ws.Range("D4").Value = 89877
$ws.Range("E4").Value = 749144
$ws.Range("F4").Value = 15042
$ws.Range("G4").Value = 134
$ws.Range("H4").Value = 50370

# Row 9 - Reino Unido
$ws.Range("B9").Value = 143464
$ws.Range("C9").Value = 5386
$ws.Range("E9").Value = 123614
$ws.Range("G9").Value = 768
$ws.Range("H9").Value = 19506

# Row 17 - Paises Bajos
$ws.Range("F17").Value = 963

# Row 25 - Austria
$ws.Range("F25").Value = 156

# Row 80 - Afganistan
$ws.Range("B80").Value = 1351
$ws.Range("C80").Value = 72
$ws.Range("E80").Value = 1120

# Row 84 - Bulgaria
$ws.Range("B84").Value = 1188
$ws.Range("C84").Value = 91
$ws.Range("E84").Value = 941
$ws.Range("G84").Value = 2
$ws.Range("H84").Value = 54

# Row 111 - Sri Lanka
$ws.Range("B111").Value = 416
$ws.Range("C111").Value = 48
$ws.Range("E111").Value = 300
